$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L1").Value = "big_decimal_value"
$ws.Range("L1").Font.Bold = $true
$ws.Range("L1").Font.Color = 0

$ws.Range("L4").Value = 93249932943949
$ws.Range("L4").NumberFormat = "0.00"

$ws.Range("L3").NumberFormat = "@"
$ws.Range("L3").Value = "3294832483943920"
$ws.Range("L3").Font.Bold = $true
$ws.Range("L3").Font.Color = 0

$ws.Range("L2").Value = -323344343.22349
$ws.Range("L2").NumberFormat = "0.0000001"
$ws.Range("L2").NumberFormat = "0.0000002"
$ws.Range("L2").NumberFormat = "0.000000"

Write-Output "done"
